$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$defaultStyleRange = $ws.Range("F2")

$ws.Range("D2").Value = "'309.04"
$ws.Range("D2").Style = $defaultStyleRange.Style
$ws.Range("E2").Value = "'0.46%"
$ws.Range("E2").Style = $defaultStyleRange.Style
$ws.Range("D3").Value = "'40.95"
$ws.Range("D3").Style = $defaultStyleRange.Style
$ws.Range("E3").Value = "'-0.13%"
$ws.Range("E3").Style = $defaultStyleRange.Style
$ws.Range("E4").Value = "'1.51%"
$ws.Range("E4").Style = $defaultStyleRange.Style
$ws.Range("E5").Value = "'0.26%"
$ws.Range("E5").Style = $defaultStyleRange.Style
$ws.Range("D6").Value = "'4.288"
$ws.Range("D6").Style = $defaultStyleRange.Style
$ws.Range("E6").Value = "'0.41%"
$ws.Range("E6").Style = $defaultStyleRange.Style
$ws.Range("D7").Value = "'1.606"
$ws.Range("D7").Style = $defaultStyleRange.Style
$ws.Range("E7").Value = "'0.29%"
$ws.Range("E7").Style = $defaultStyleRange.Style
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.463"
$ws.Range("D8").Style = $defaultStyleRange.Style
$ws.Range("E8").Value = "'1.11%"
$ws.Range("E8").Style = $defaultStyleRange.Style
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9097"
$ws.Range("D9").Style = $defaultStyleRange.Style
$ws.Range("E10").Value = "'24.71%"
$ws.Range("E10").Style = $defaultStyleRange.Style
$ws.Range("D11").Value = "'0.1805"
$ws.Range("D11").Style = $defaultStyleRange.Style
$ws.Range("E11").Value = "'2.36%"
$ws.Range("E11").Style = $defaultStyleRange.Style
$ws.Range("D12").Value = "'0.09138"
$ws.Range("D12").Style = $defaultStyleRange.Style
$ws.Range("E12").Value = "'0.61%"
$ws.Range("E12").Style = $defaultStyleRange.Style
$ws.Range("D13").Value = "'0.04331"
$ws.Range("D13").Style = $defaultStyleRange.Style
$ws.Range("E13").Value = "'-0.41%"
$ws.Range("E13").Style = $defaultStyleRange.Style
$ws.Range("E14").Value = "'-0.54%"
$ws.Range("E14").Style = $defaultStyleRange.Style
$ws.Range("D15").Value = "'0.001247"
$ws.Range("D15").Style = $defaultStyleRange.Style
$ws.Range("E15").Value = "'-0.68%"
$ws.Range("E15").Style = $defaultStyleRange.Style
$ws.Range("D16").Value = "'0.005663"
$ws.Range("D16").Style = $defaultStyleRange.Style
$ws.Range("E16").Value = "'-3.70%"
$ws.Range("E16").Style = $defaultStyleRange.Style
$ws.Range("D17").Value = "'3.348"
$ws.Range("D17").Style = $defaultStyleRange.Style
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("E17").Style = $defaultStyleRange.Style
$ws.Range("E18").Value = "'1.49%"
$ws.Range("E18").Style = $defaultStyleRange.Style
$ws.Range("D19").Value = "'6.988"
$ws.Range("D19").Style = $defaultStyleRange.Style
$ws.Range("E19").Value = "'2.56%"
$ws.Range("E19").Style = $defaultStyleRange.Style
$ws.Range("D20").Value = "'0.1393"
$ws.Range("D20").Style = $defaultStyleRange.Style
$ws.Range("E20").Value = "'2.60%"
$ws.Range("E20").Style = $defaultStyleRange.Style
$ws.Range("D21").Value = "'0.2739"
$ws.Range("D21").Style = $defaultStyleRange.Style
$ws.Range("E21").Value = "'0.42%"
$ws.Range("E21").Style = $defaultStyleRange.Style
$ws.Range("D22").Value = "'0.04039"
$ws.Range("D22").Style = $defaultStyleRange.Style
$ws.Range("E22").Value = "'-3.41%"
$ws.Range("E22").Style = $defaultStyleRange.Style
$ws.Range("E23").Value = "'3.50%"
$ws.Range("E23").Style = $defaultStyleRange.Style
$ws.Range("D24").Value = "'0.004054"
$ws.Range("D24").Style = $defaultStyleRange.Style
$ws.Range("E24").Value = "'-0.24%"
$ws.Range("E24").Style = $defaultStyleRange.Style
$ws.Range("E25").Value = "'-2.41%"
$ws.Range("E25").Style = $defaultStyleRange.Style
$ws.Range("E26").Value = "'24.48%"
$ws.Range("E26").Style = $defaultStyleRange.Style
$ws.Range("D38").Value = "'0.02419"
$ws.Range("D38").Style = $defaultStyleRange.Style
$ws.Range("E38").Value = "'0.93%"
$ws.Range("E38").Style = $defaultStyleRange.Style
$ws.Range("D39").Value = "'0.05232"
$ws.Range("D39").Style = $defaultStyleRange.Style
$ws.Range("E39").Value = "'0.94%"
$ws.Range("E39").Style = $defaultStyleRange.Style
$ws.Range("D40").Value = "'0.007836"
$ws.Range("D40").Style = $defaultStyleRange.Style
$ws.Range("E40").Value = "'0.33%"
$ws.Range("E40").Style = $defaultStyleRange.Style
$ws.Range("D41").Value = "'0.1301"
$ws.Range("D41").Style = $defaultStyleRange.Style
$ws.Range("E41").Value = "'-0.12%"
$ws.Range("E41").Style = $defaultStyleRange.Style
$ws.Range("D42").Value = "'0.006810"
$ws.Range("D42").Style = $defaultStyleRange.Style
$ws.Range("E42").Value = "'-4.00%"
$ws.Range("E42").Style = $defaultStyleRange.Style
$ws.Range("D43").Value = "'0.001862"
$ws.Range("D43").Style = $defaultStyleRange.Style
$ws.Range("E43").Value = "'-3.09%"
$ws.Range("E43").Style = $defaultStyleRange.Style
$ws.Range("E44").Value = "'-0.75%"
$ws.Range("E44").Style = $defaultStyleRange.Style
$ws.Range("D45").Value = "'0.3348"
$ws.Range("D45").Style = $defaultStyleRange.Style
$ws.Range("E45").Value = "'-0.12%"
$ws.Range("E45").Style = $defaultStyleRange.Style
$ws.Range("D46").Value = "'0.00006870"
$ws.Range("D46").Style = $defaultStyleRange.Style
$ws.Range("E46").Value = "'8.05%"
$ws.Range("E46").Style = $defaultStyleRange.Style
$ws.Range("E47").Value = "'-0.10%"
$ws.Range("E47").Style = $defaultStyleRange.Style
$ws.Range("D48").Value = "'0.1426"
$ws.Range("D48").Style = $defaultStyleRange.Style
$ws.Range("E48").Value = "'2,335.85%"
$ws.Range("E48").Style = $defaultStyleRange.Style
$ws.Range("E49").Value = "'-31.90%"
$ws.Range("E49").Style = $defaultStyleRange.Style
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = $defaultStyleRange.Style
$ws.Range("E50").Value = "'-0.10%"
$ws.Range("E50").Style = $defaultStyleRange.Style
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = $defaultStyleRange.Style
$ws.Range("E51").Value = "'-0.10%"
$ws.Range("E51").Style = $defaultStyleRange.Style
